# add html5 & git
# Populate the two previously-empty rows (4 and 5) in the tools list on
# Sheet1 with the new "HTML5" and "git" entries, widen column C so the
# (now much longer) URLs are visible, and leave the selection where the
# author's Excel session ended up (C5).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4: HTML5
$ws.Range("B4").Value = "HTML5"
$ws.Range("C4").Value = "http://e.jikexueyuan.com/html5.html?hmsr=baidu_sem_html5_dy_5"

# Row 5: git
$ws.Range("B5").Value = "git"
$ws.Range("C5").Value = "http://www.liaoxuefeng.com/wiki/0013739516305929606dd18361248578c67b8067c8c017b000"

# Column C needs to be wide enough for the new long URLs (target raw OOXML
# width 80.125). The engine snaps ColumnWidth to a 1/7-character pixel grid
# (pixels = round(ColumnWidth*7); storedWidth = (pixels+5)/7), so feed it
# the ColumnWidth value whose grid-snapped result lands closest to 80.125.
$ws.Columns("C").ColumnWidth = 79.42857142857143

# Leave the active selection on C5, matching the saved sheet view.
[void]$ws.Range("C5").Select()
